$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delay before executing the remaining script actions (per commit message:
# "Adding the delay wwhile executing scripts")
Start-Sleep -Milliseconds 500

# --- Update row 2 candidate data ---
$ws.Range("A2").Value = "Ashwini.shashikiran"
$ws.Range("F2").Value = "female"
$ws.Range("G2").Value = "George.Thomps5@gmail.com"
$ws.Range("H2").Value = 9741545882
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 60

# --- Data validations ---
$ws.Range("I7").Validation.Add(3, 1, 1, '"Fresher,Experienced"')
$ws.Range("I5").Validation.Add(3, 1, 1, '"Experienced,Fresher"')
$ws.Range("I9").Validation.Add(3, 1, 1, '"Experienced,Fresher"')
$ws.Range("I2:I3").Validation.Add(3, 1, 1, '"Experienced,Fresher"')
$ws.Range("F2:F3").Validation.Add(3, 1, 1, '"female,male"')
$ws.Range("D2:D8").Validation.Add(3, 1, 1, '"IT,BPO,PST,SSS-Shared Services,SHILOH,GC-IT,DIGITAL"')

# --- Selection / view: move from G2 (scrolled to D1) back to B2 at default scroll ---
$ws.Range("B2").Select()
